$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 7 down to row 8 so the new row matches the
# existing table row styling (fill, number formats, etc.)
$ws.Range("A7").Copy($ws.Range("A8"))
$ws.Range("B7").Copy($ws.Range("B8"))
$ws.Range("C7").Copy($ws.Range("C8"))
$ws.Range("D7").Copy($ws.Range("D8"))
$ws.Range("E7").Copy($ws.Range("E8"))
$ws.Range("F7").Copy($ws.Range("F8"))
$ws.Range("G7").Copy($ws.Range("G8"))

# New data row: geothermal
$ws.Range("A8").Value = "geothermal"
$ws.Range("B8").Value = 5000000
$ws.Range("C8").Formula = "=0.3"
$ws.Range("D8").Formula = "=(1-Table1[[#This Row],[Direct Pay ITC]])*Table1[[#This Row],[$/MW Gross Capital Cost]]"
$ws.Range("E8").Value = 0.029
$ws.Range("F8").Value = 20
$ws.Range("G8").Formula = "=PMT(Table1[[#This Row],[Annual Rate]],Table1[[#This Row],[Term]],Table1[[#This Row],[Net Capital Cost]])"

# Grow Table1 to include the new row
$ws.ListObjects("Table1").Resize($ws.Range("A1:I8"))

$ws.Range("B9").Select()

$wb.Application.Calculate()
